$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Opp")
$ws.Activate()

$ws.Range("A2").Value = "TestAutomation1"
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = "Needs Analysis"

$ws.Range("A3").Value = "TestAutomation2"
$ws.Range("B3").Value = 21
$ws.Range("C3").Value = "Needs Analysis"

$ws.Range("A4").Value = "TestAutomation3"
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = "Needs Analysis"

$ws.Range("A5").Value = "TestAutomation4"
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = "Needs Analysis"

$ws.Range("A6").Value = "TestAutomation5"
$ws.Range("B6").Value = 24
$ws.Range("C6").Value = "Needs Analysis"

$ws.Range("A7").Value = "TestAutomation6"
$ws.Range("B7").Value = 25
$ws.Range("C7").Value = "Needs Analysis"

$ws.Range("A8").Value = "TestAutomation7"
$ws.Range("B8").Value = 26
$ws.Range("C8").Value = "Needs Analysis"

$ws.Range("A9").Value = "TestAutomation8"
$ws.Range("B9").Value = 27
$ws.Range("C9").Value = "Needs Analysis"

$ws.Range("B3").Select()
